$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '26.109.79'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").Value = "'" + '1.597.66'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'" + '212.02'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = "'" + '0.485'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.65%  '
$ws.Range("D8").Value = "'" + '0.248'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.20%  '
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("D10").Value = "'" + '17.90'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("D11").Value = "'" + '0.0820'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.81%  '
$ws.Range("D12").Value = "'" + '1.821.68'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.12%  '
$ws.Range("D13").Value = "'" + '1.605.81'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").Value = "'" + '26.086.97'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.18%  '
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = "'" + '0.0₃0721'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = "'" + '204.66'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +10.06%  '
$ws.Range("E21").Value = '  +2.20%  '
$ws.Range("D22").Value = "'" + '9.30'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  +11.08%  '
$ws.Range("D25").Value = "'" + '141.94'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.64%  '
$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  -3.15%  '
$ws.Range("D28").Value = "'" + '15.22'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.42%  '
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("E35").Value = '  +2.58%  '
$ws.Range("D36").Value = "'" + '0.0164'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.62%  '
$ws.Range("D37").Value = "'" + '1.106.39'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").Value = "'" + '0.778'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("E42").Value = '  +1.78%  '
$ws.Range("D43").Value = "'" + '1.734.45'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.13%  '
$ws.Range("D44").Value = "'" + '92.56'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("E45").Value = '  +0.35%  '
$ws.Range("E46").Value = '  +6.41%  '
$ws.Range("D47").Value = "'" + '0.0₆0102'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.24%  '
$ws.Range("D48").Value = "'" + '53.32'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.99%  '
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("E51").Value = '  +0.05%  '
